$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write an exact text value to a cell without Excel re-typing it
# (numbers, percentages) as a number - force text via a leading apostrophe
# then strip the resulting quote-prefix style so the cell stays "Normal".
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "305.13"
Set-TextCell $ws "E2" "0.67%"

Set-TextCell $ws "D3" "35.80"
Set-TextCell $ws "E3" "1.49%"

Set-TextCell $ws "D4" "5.081"
Set-TextCell $ws "E4" "-0.03%"

Set-TextCell $ws "D5" "0.08078"
Set-TextCell $ws "E5" "0.84%"

Set-TextCell $ws "D6" "1.920"
Set-TextCell $ws "E6" "-0.84%"

Set-TextCell $ws "D7" "4.190"
Set-TextCell $ws "E7" "3.18%"

Set-TextCell $ws "D8" "7.756"
Set-TextCell $ws "E8" "-0.12%"

Set-TextCell $ws "D9" "0.9265"
Set-TextCell $ws "E9" "0.24%"

Set-TextCell $ws "D10" "0.1370"
Set-TextCell $ws "E10" "11.36%"

Set-TextCell $ws "D11" "0.1892"
Set-TextCell $ws "E11" "2.13%"

Set-TextCell $ws "D12" "0.09209"
Set-TextCell $ws "E12" "-5.05%"

Set-TextCell $ws "D13" "0.03421"
Set-TextCell $ws "E13" "-5.68%"

Set-TextCell $ws "D14" "0.09813"
Set-TextCell $ws "E14" "-0.47%"

Set-TextCell $ws "D15" "0.001450"
Set-TextCell $ws "E15" "4.35%"

Set-TextCell $ws "D16" "0.005793"
Set-TextCell $ws "E16" "-0.26%"

Set-TextCell $ws "D17" "3.621"
Set-TextCell $ws "E17" "3.42%"

Set-TextCell $ws "D18" "3.015"
Set-TextCell $ws "E18" "1.90%"

Set-TextCell $ws "D19" "0.3456"
Set-TextCell $ws "E19" "1.47%"

Set-TextCell $ws "D20" "0.1334"
Set-TextCell $ws "E20" "1.79%"

Set-TextCell $ws "D21" "4.927"
Set-TextCell $ws "E21" "-2.46%"

Set-TextCell $ws "D22" "0.2443"
Set-TextCell $ws "E22" "-0.98%"

Set-TextCell $ws "D23" "0.04438"
Set-TextCell $ws "E23" "-1.98%"

Set-TextCell $ws "E24" "0.15%"

Set-TextCell $ws "D25" "0.004815"
Set-TextCell $ws "E25" "-0.42%"

Set-TextCell $ws "D26" "0.0001301"
Set-TextCell $ws "E26" "3.97%"

Set-TextCell $ws "D27" "0.0003131"
Set-TextCell $ws "E27" "4.23%"

Set-TextCell $ws "D39" "0.02017"
Set-TextCell $ws "E39" "4.79%"

Set-TextCell $ws "D40" "0.04922"
Set-TextCell $ws "E40" "4.55%"

Set-TextCell $ws "D41" "0.007628"
Set-TextCell $ws "E41" "0.99%"

Set-TextCell $ws "D42" "0.01030"
Set-TextCell $ws "E42" "5.99%"

Set-TextCell $ws "D43" "0.1373"
Set-TextCell $ws "E43" "3.32%"

Set-TextCell $ws "D44" "0.002101"
Set-TextCell $ws "E44" "-0.51%"

Set-TextCell $ws "D45" "0.01105"
Set-TextCell $ws "E45" "9.95%"

Set-TextCell $ws "D46" "0.00006417"
Set-TextCell $ws "E46" "2.17%"

Set-TextCell $ws "D49" "0.001192"
Set-TextCell $ws "E49" "-19.97%"

Set-TextCell $ws "D50" "0.00002102"

Set-TextCell $ws "D51" "0.0002002"
